$d = $word.ActiveDocument

$replacements = @(
    @("96÷4=24, 0", "50÷7=7, 1"),
    @("89÷9=9, 8", "55÷9=6, 1"),
    @("17÷4=4, 1", "93÷7=13, 2"),
    @("62÷6=10, 2", "19÷8=2, 3"),
    @("68÷9=7, 5", "10÷6=1, 4"),
    @("84÷7=12, 0", "50÷5=10, 0"),
    @("42÷5=8, 2", "73÷5=14, 3"),
    @("29÷4=7, 1", "33÷5=6, 3"),
    @("18÷9=2, 0", "17÷8=2, 1"),
    @("90÷8=11, 2", "37÷2=18, 1"),
    @("21÷2=10, 1", "55÷7=7, 6"),
    @("33÷3=11, 0", "20÷9=2, 2"),
    @("56÷3=18, 2", "45÷7=6, 3"),
    @("40÷9=4, 4", "28÷9=3, 1"),
    @("13÷8=1, 5", "49÷5=9, 4"),
    @("88÷2=44, 0", "76÷9=8, 4"),
    @("75÷8=9, 3", "71÷4=17, 3"),
    @("87÷6=14, 3", "71÷3=23, 2"),
    @("66÷4=16, 2", "75÷9=8, 3"),
    @("42÷6=7, 0", "30÷3=10, 0"),
    @("64÷5=12, 4", "51÷9=5, 6"),
    @("53÷2=26, 1", "73÷7=10, 3"),
    @("25÷8=3, 1", "49÷8=6, 1"),
    @("45÷3=15, 0", "62÷4=15, 2"),
    @("74÷8=9, 2", "54÷5=10, 4")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
